$newValues = @(
    "35×49=1715",
    "57×23=1311",
    "72×43=3096",
    "49×99=4851",
    "94×23=2162",
    "42×15=630",
    "52×95=4940",
    "53×71=3763",
    "54×55=2970",
    "38×34=1292",
    "92×72=6624",
    "16×76=1216",
    "87×53=4611",
    "41×33=1353",
    "85×48=4080",
    "72×75=5400",
    "37×12=444",
    "30×92=2760",
    "15×88=1320",
    "74×44=3256",
    "78×23=1794",
    "49×23=1127",
    "86×89=7654",
    "18×64=1152",
    "78×16=1248",
    "40×51=2040",
    "52×100=5200",
    "44×74=3256",
    "61×68=4148",
    "71×82=5822",
    "80×75=6000",
    "17×52=884",
    "85×56=4760",
    "18×71=1278",
    "40×42=1680",
    "11×16=176",
    "16×20=320",
    "14×19=266",
    "100×49=4900",
    "28×56=1568",
    "20×54=1080",
    "57×30=1710",
    "43×48=2064",
    "14×90=1260",
    "41×60=2460",
    "20×17=340",
    "30×43=1290",
    "39×76=2964",
    "55×12=660",
    "64×51=3264",
    "25×43=1075",
    "38×53=2014",
    "25×54=1350",
    "47×12=564",
    "48×72=3456",
    "64×56=3584",
    "47×91=4277",
    "22×29=638",
    "13×80=1040",
    "60×93=5580",
    "44×22=968",
    "25×52=1300",
    "34×30=1020",
    "35×53=1855",
    "80×84=6720",
    "77×52=4004",
    "49×54=2646",
    "37×62=2294",
    "75×58=4350",
    "14×71=994",
    "84×28=2352",
    "20×16=320",
    "73×45=3285",
    "74×88=6512",
    "34×11=374",
    "96×45=4320",
    "36×75=2700",
    "28×36=1008",
    "40×49=1960",
    "57×93=5301",
    "62×18=1116",
    "44×43=1892",
    "70×85=5950",
    "73×93=6789",
    "89×56=4984",
    "76×21=1596",
    "31×72=2232",
    "14×63=882",
    "33×17=561",
    "41×10=410",
    "68×71=4828",
    "74×89=6586",
    "10×33=330",
    "74×84=6216",
    "94×36=3384",
    "10×82=820",
    "18×67=1206",
    "43×54=2322",
    "45×23=1035",
    "77×82=6314"
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host ("Updated " + $idx + " cells")
